# Modify command style for both earnings and payment
# This script edits three sequence-diagram textboxes on slide 1, replacing the
# verbose "idx/1 amt/200 m/8 y/2018" style command-argument text with the
# shorter "1 200 8 2018" style, and shrinking/repositioning the textboxes to
# match their new (single/double-line) content.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$LDQ = [char]8220   # “
$RDQ = [char]8221   # ”

# ---------------------------------------------------------------------------
# Shape "TextBox 25" (cNvPr id=26) -> execute( "paid idx/1 amt/200 m/8 y/2018")
# ---------------------------------------------------------------------------
$shp1 = $s.Shapes.Item(11)
$tr1 = $shp1.TextFrame.TextRange
$len1 = $tr1.Text.Length
$full1 = $tr1.Characters(1, $len1)
$full1.Text = "execute( " + $LDQ + "paid 1 200 8 2018" + $RDQ + ")"

$shp1.Left = 5.736299276598425
$shp1.Top = 76.26110458420473
$shp1.Width = 148.0408706667323
$shp1.Height = 13.328897953795275

# ---------------------------------------------------------------------------
# Shape "TextBox 79" (cNvPr id=80) -> parseCommand("paid idx/1 amt/200 m/8 y/2018")
# ---------------------------------------------------------------------------
$shp2 = $s.Shapes.Item(18)
$tr2 = $shp2.TextFrame.TextRange
$len2 = $tr2.Text.Length
# Keep the leading "parseCommand" run (12 chars) untouched; rewrite the rest.
$rest2 = $tr2.Characters(13, $len2 - 12)
$rest2.Text = "(" + $LDQ + "paid 1 200 8 2018" + $RDQ + ")"

$shp2.Height = 13.328897953795275

# ---------------------------------------------------------------------------
# Shape "TextBox 60" (cNvPr id=61) -> Parse("idx/1 amt/200 m/8 y/2018")
# ---------------------------------------------------------------------------
$shp3 = $s.Shapes.Item(46)
$tr3 = $shp3.TextFrame.TextRange
$len3 = $tr3.Text.Length
# Replace the trailing "/1 amt/200 m/8 y/2018”)" with just ”)
$tail3 = $tr3.Characters(11, $len3 - 10)
$tail3.Text = $RDQ + ")"
# Replace the middle "idx" run with the new digits.
$mid3 = $tr3.Characters(8, 3)
$mid3.Text = "1 200 8 2018"
# Merge the leading "Parse(“" run together with the digits we just inserted.
$len3b = $tr3.Text.Length
$prefix3 = $tr3.Characters(1, $len3b - 2)
$prefix3.Text = "Parse(" + $LDQ + "1 200 8 2018"

$shp3.Height = 36.3515758521496
